# Update the "Price" (column D) values for the symbol list snapshot.
# Values are stored as text (inline strings), so we force a Text number
# format before writing and then restore the default style, keeping the
# cells looking exactly like every other untouched text cell on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    2 = "268.13"
    3 = "22.98"
    4 = "6.311"
    5 = "0.06186"
    6 = "3.590"
    7 = "6.691"
    8 = "1.387"
    9 = "0.8314"
    11 = "0.1605"
    12 = "0.08243"
    13 = "0.03420"
    14 = "0.03246"
    15 = "0.09277"
    16 = "3.901"
    17 = "0.001723"
    18 = "0.04844"
    19 = "0.006283"
    20 = "0.005375"
    21 = "0.001090"
    22 = "0.0001500"
    23 = "3.772"
    24 = "2.366"
    27 = "0.0002683"
    40 = "0.04674"
    41 = "0.006971"
    43 = "0.003351"
    44 = "0.01220"
    45 = "0.00006239"
    47 = "0.7003"
    48 = "0.1769"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)   # column D = Price
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}
